{"js": "// Replace the date heading and each math problem's operands with the new\n// values from the target revision. Replacements are matched strictly in\n// document order so a value that is both an old text (e.g. \"927\\u00d78=\")\n// and a new text elsewhere never gets double-replaced: every search() is\n// issued against the *original* content, and only after all the matching\n// ranges have been located do we write the new text into them.\nconst replacements = [\n  { oldText: \"2026-01-10 Saturday\", newText: \"2026-01-11 Sunday\" },\n  { oldText: \"673\u00d74=\", newText: \"107\u00d75=\" },\n  { oldText: \"551\u00d72=\", newText: \"927\u00d78=\" },\n  { oldText: \"682\u00d74=\", newText: \"427\u00d73=\" },\n  { oldText: \"680\u00d73=\", newText: \"942\u00d72=\" },\n  { oldText: \"271\u00d74=\", newText: \"983\u00d74=\" },\n  { oldText: \"246\u00d79=\", newText: \"930\u00d78=\" },\n  { oldText: \"471\u00d72=\", newText: \"511\u00d75=\" },\n  { oldText: \"288\u00d73=\", newText: \"228\u00d76=\" },\n  { oldText: \"559\u00d74=\", newText: \"416\u00d76=\" },\n  { oldText: \"503\u00d72=\", newText: \"129\u00d74=\" },\n  { oldText: \"832\u00d72=\", newText: \"380\u00d73=\" },\n  { oldText: \"927\u00d78=\", newText: \"230\u00d74=\" },\n  { oldText: \"850\u00d73=\", newText: \"347\u00d74=\" },\n  { oldText: \"564\u00d72=\", newText: \"804\u00d73=\" },\n  { oldText: \"282\u00d75=\", newText: \"616\u00d74=\" },\n  { oldText: \"848\u00d79=\", newText: \"282\u00d72=\" },\n  { oldText: \"462\u00d73=\", newText: \"290\u00d73=\" },\n  { oldText: \"615\u00d76=\", newText: \"378\u00d78=\" },\n  { oldText: \"541\u00d77=\", newText: \"780\u00d76=\" },\n  { oldText: \"837\u00d75=\", newText: \"846\u00d75=\" },\n  { oldText: \"358\u00d73=\", newText: \"405\u00d78=\" },\n  { oldText: \"571\u00d79=\", newText: \"542\u00d79=\" },\n  { oldText: \"564\u00d75=\", newText: \"432\u00d79=\" },\n  { oldText: \"817\u00d73=\", newText: \"704\u00d74=\" },\n  { oldText: \"827\u00d73=\", newText: \"276\u00d73=\" },\n];\n\nconst body = context.document.body;\n\n// First, locate every occurrence (one each, since all the source strings\n// are unique in this document) while the body still holds the original text.\nconst searchResults = replacements.map((r) =>\n  body.search(r.oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((results) => results.load(\"items\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const results = searchResults[i];\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${replacements[i].oldText}`);\n  }\n  results.items[0].insertText(replacements[i].newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n", "ps1": "# Update the date heading and each math-problem cell with the values from\n# the target revision. Cells are addressed by their fixed (row, column)\n# position in the single table rather than by searching for the old text,\n# because one of the new values (\"927\" + [char]0x00D7 + \"8=\") is identical\n# to another cell's *original* value elsewhere in the table; a content-based\n# Find/ReplaceAll pass run repeatedly would otherwise clobber that unrelated\n# cell once its own replacement text happened to match a still-pending search\n# string. Addressing by cell position sidesteps that entirely.\n\n$d = $word.ActiveDocument\n\n$dateRange = $d.Paragraphs(1).Range\n$dateCurrent = $dateRange.Text.TrimEnd([char]13)\nif ($dateCurrent -ne '2026-01-10 Saturday') {\n    throw \"Unexpected date paragraph text: '$dateCurrent'\"\n}\n$dateRange.Text = '2026-01-11 Sunday'\n\n$tbl = $d.Tables(1)\n\n$cell = $tbl.Cell(1, 1)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '673\u00d74=') {\n    throw \"Unexpected text in cell (1,1): '$cellCurrent'\"\n}\n$cell.Range.Text = '107\u00d75='\n\n$cell = $tbl.Cell(1, 2)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '551\u00d72=') {\n    throw \"Unexpected text in cell (1,2): '$cellCurrent'\"\n}\n$cell.Range.Text = '927\u00d78='\n\n$cell = $tbl.Cell(1, 3)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '682\u00d74=') {\n    throw \"Unexpected text in cell (1,3): '$cellCurrent'\"\n}\n$cell.Range.Text = '427\u00d73='\n\n$cell = $tbl.Cell(1, 4)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '680\u00d73=') {\n    throw \"Unexpected text in cell (1,4): '$cellCurrent'\"\n}\n$cell.Range.Text = '942\u00d72='\n\n$cell = $tbl.Cell(1, 5)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '271\u00d74=') {\n    throw \"Unexpected text in cell (1,5): '$cellCurrent'\"\n}\n$cell.Range.Text = '983\u00d74='\n\n$cell = $tbl.Cell(5, 1)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '246\u00d79=') {\n    throw \"Unexpected text in cell (5,1): '$cellCurrent'\"\n}\n$cell.Range.Text = '930\u00d78='\n\n$cell = $tbl.Cell(5, 2)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '471\u00d72=') {\n    throw \"Unexpected text in cell (5,2): '$cellCurrent'\"\n}\n$cell.Range.Text = '511\u00d75='\n\n$cell = $tbl.Cell(5, 3)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '288\u00d73=') {\n    throw \"Unexpected text in cell (5,3): '$cellCurrent'\"\n}\n$cell.Range.Text = '228\u00d76='\n\n$cell = $tbl.Cell(5, 4)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '559\u00d74=') {\n    throw \"Unexpected text in cell (5,4): '$cellCurrent'\"\n}\n$cell.Range.Text = '416\u00d76='\n\n$cell = $tbl.Cell(5, 5)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '503\u00d72=') {\n    throw \"Unexpected text in cell (5,5): '$cellCurrent'\"\n}\n$cell.Range.Text = '129\u00d74='\n\n$cell = $tbl.Cell(10, 1)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '832\u00d72=') {\n    throw \"Unexpected text in cell (10,1): '$cellCurrent'\"\n}\n$cell.Range.Text = '380\u00d73='\n\n$cell = $tbl.Cell(10, 2)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '927\u00d78=') {\n    throw \"Unexpected text in cell (10,2): '$cellCurrent'\"\n}\n$cell.Range.Text = '230\u00d74='\n\n$cell = $tbl.Cell(10, 3)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '850\u00d73=') {\n    throw \"Unexpected text in cell (10,3): '$cellCurrent'\"\n}\n$cell.Range.Text = '347\u00d74='\n\n$cell = $tbl.Cell(10, 4)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '564\u00d72=') {\n    throw \"Unexpected text in cell (10,4): '$cellCurrent'\"\n}\n$cell.Range.Text = '804\u00d73='\n\n$cell = $tbl.Cell(10, 5)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '282\u00d75=') {\n    throw \"Unexpected text in cell (10,5): '$cellCurrent'\"\n}\n$cell.Range.Text = '616\u00d74='\n\n$cell = $tbl.Cell(15, 1)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '848\u00d79=') {\n    throw \"Unexpected text in cell (15,1): '$cellCurrent'\"\n}\n$cell.Range.Text = '282\u00d72='\n\n$cell = $tbl.Cell(15, 2)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '462\u00d73=') {\n    throw \"Unexpected text in cell (15,2): '$cellCurrent'\"\n}\n$cell.Range.Text = '290\u00d73='\n\n$cell = $tbl.Cell(15, 3)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '615\u00d76=') {\n    throw \"Unexpected text in cell (15,3): '$cellCurrent'\"\n}\n$cell.Range.Text = '378\u00d78='\n\n$cell = $tbl.Cell(15, 4)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '541\u00d77=') {\n    throw \"Unexpected text in cell (15,4): '$cellCurrent'\"\n}\n$cell.Range.Text = '780\u00d76='\n\n$cell = $tbl.Cell(15, 5)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '837\u00d75=') {\n    throw \"Unexpected text in cell (15,5): '$cellCurrent'\"\n}\n$cell.Range.Text = '846\u00d75='\n\n$cell = $tbl.Cell(20, 1)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '358\u00d73=') {\n    throw \"Unexpected text in cell (20,1): '$cellCurrent'\"\n}\n$cell.Range.Text = '405\u00d78='\n\n$cell = $tbl.Cell(20, 2)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '571\u00d79=') {\n    throw \"Unexpected text in cell (20,2): '$cellCurrent'\"\n}\n$cell.Range.Text = '542\u00d79='\n\n$cell = $tbl.Cell(20, 3)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '564\u00d75=') {\n    throw \"Unexpected text in cell (20,3): '$cellCurrent'\"\n}\n$cell.Range.Text = '432\u00d79='\n\n$cell = $tbl.Cell(20, 4)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '817\u00d73=') {\n    throw \"Unexpected text in cell (20,4): '$cellCurrent'\"\n}\n$cell.Range.Text = '704\u00d74='\n\n$cell = $tbl.Cell(20, 5)\n$cellRange = $cell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellCurrent = $cellRange.Text.TrimEnd([char]13)\nif ($cellCurrent -ne '827\u00d73=') {\n    throw \"Unexpected text in cell (20,5): '$cellCurrent'\"\n}\n$cell.Range.Text = '276\u00d73='\n\n"}
